$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column (H1), cloning the format of the
# existing header cell (G1) so it picks up the same bold/border/
# center-top style used by the rest of the header row.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the corresponding data cell in row 2.
$ws.Range("H2").Value = 1
